$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 241, pushing the existing data (old rows
# 241-317) down to rows 242-318. This also grows the used range from
# A1:R317 to A1:R318 automatically.
$ws.Rows(241).Insert()

# Populate the newly inserted row 241 with the new record (same
# market/category metadata as its neighbour, but its own date, volume,
# price and origin values).
$ws.Range("A241").Value = 4
$ws.Range("B241").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C241").Value = "Los Lagos"
$ws.Range("D241").Value = 44627
$ws.Range("E241").Value = 10
$ws.Range("F241").Value = 100114013
$ws.Range("G241").Value = "Zanahoria"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 250
$ws.Range("K241").Value = 12000
$ws.Range("L241").Value = 12000
$ws.Range("M241").Value = 12000
$ws.Range("N241").Value = "`$/saco 20 kilos"
$ws.Range("O241").Value = "Chillán"
$ws.Range("P241").Value = 600
$ws.Range("Q241").Value = 20
$ws.Range("R241").Value = "Hortaliza"
